# Updating n_total_regvars to be length_unique_varnames in the
# modelsummary_reg_default table: this adds two additional "NA" rows
# (one per newly-unique model variable row) to the regression summary
# table, inserted directly above the "R2" row (i.e. right after the
# "Mean" row).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the "R2" row -- the two new rows of "NA" placeholders belong
# directly above it (right after the "Mean" row).
$anchorRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $label = $t.Rows.Item($i).Cells.Item(1).Range.Text
    if ($label -like "R2*") {
        $anchorRow = $t.Rows.Item($i)
        break
    }
}

if ($anchorRow -eq $null) {
    $anchorRow = $t.Rows.Item($t.Rows.Count)
}

$numColumns = $t.Columns.Count

for ($r = 0; $r -lt 2; $r++) {
    $newRow = $t.Rows.Add($anchorRow)
    for ($c = 1; $c -le $newRow.Cells.Count; $c++) {
        $newRow.Cells.Item($c).Range.Text = "NA"
    }
}
